# "updated suzie excel file." - Suzie's Roth IRA received a $2.99 dividend
# in December, which needs to be reflected on the Yearly sheet (and which
# flows through to the All Time sheet's running totals).

$wb = $excel.ActiveWorkbook
$wsYearly = $wb.Worksheets.Item("Yearly")
$wsAllTime = $wb.Worksheets.Item("All Time")

# December is row 14 on the Yearly sheet: B=month#, C=month name,
# D=Taxable Account, E=401K, F=Suzie's Roth IRA, G=Grand Total.
$wsYearly.Range("F14").Value = 2.99
$wsYearly.Range("G14").Formula = "=SUM(D14:F14)"

# Leave the selection/active-sheet state the way the author left it when
# they saved after making the edit on the Yearly sheet.
$wsAllTime.Range("M41").Select() | Out-Null
$wsYearly.Activate() | Out-Null
$wsYearly.Range("I5").Select() | Out-Null

$wb.Save()
